$wb = $excel.ActiveWorkbook

$wsStartPrice = $wb.Worksheets.Item("start_price")
$wsLinear = $wb.Worksheets.Item("Linear")
$wsNonLinear = $wb.Worksheets.Item("NonLinear")

# start_price sheet
$wsStartPrice.Range("A2").Value = 3968.50512356841

# Linear sheet
$wsLinear.Range("B2").Value = 0.03773312853246484
$wsLinear.Range("B3").Value = 0.1398872671655756
$wsLinear.Range("B4").Value = 924.0722602314153
$wsLinear.Range("B5").Value = "[1.0, 0.228907528583736, 0.09538452892358093, 0.14228286298227724, 0.09928889169104207, 0.06383026520644189, 0.2015639282132029, 0.3394807651048156, 0.18311080742058392, 0.04375830110136642, 0.07838524168038673, 0.09029172830065485, 0.031147295562521993, 0.17353070558452627, 0.3226748149391493, 0.15552612993285225, 0.010828395862666302, 0.06589093519872523, 0.057323253835243365, -0.002552218692806691]"

# NonLinear sheet
$wsNonLinear.Range("B3").Value = 0.8557046979865772
$wsNonLinear.Range("B4").Value = 0.1508203477638428
$wsNonLinear.Range("B5").Value = 0.1819105945768426
$wsNonLinear.Range("B6").Value = 951.2803638735879
$wsNonLinear.Range("B7").Value = 1.408848748468536
$wsNonLinear.Range("B8").Value = 0.02304789272140188
$wsNonLinear.Range("B9").Value = 899.9544305974064
$wsNonLinear.Range("B10").Value = "[1.0, 0.2265631528631676, 0.09276183708173537, 0.1448534695558428, 0.10153516502393309, 0.0633902165221686, 0.20044946287356, 0.3366031374618973, 0.1813720031402772, 0.04388853382947027, 0.07955141366456671, 0.09165940500508725, 0.03150290542297787, 0.1723662681634362, 0.321271681104698, 0.1550717377286061, 0.011236230688369259, 0.06622689403941319, 0.057057576155967155, -0.0031778573396040477]"
